$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 7586.5654
$ws.Range("I98").Value = 7885.6816
$ws.Range("J98").Value = 1006
$ws.Range("K98").Value = 7885.6816
$ws.Range("L98").Value = 1006
$ws.Range("M98").Value = -6387.6816
$ws.Range("N98").Value = -4002
# Row 111
$ws.Range("H111").Value = 2174
$ws.Range("I111").Value = 2096.0908
$ws.Range("K111").Value = 6288.2724
$ws.Range("M111").Value = -3221.2724
# Row 116
$ws.Range("H116").Value = 3413.2
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442
# Row 122
$ws.Range("H122").Value = 7586.5654
$ws.Range("I122").Value = 7885.6816
$ws.Range("J122").Value = 1006
$ws.Range("K122").Value = 23657.0448
$ws.Range("L122").Value = 3018
$ws.Range("M122").Value = -21207.0448
$ws.Range("N122").Value = -7918
# Row 126
$ws.Range("H126").Value = 100000
$ws.Range("J126").Value = 100000
$ws.Range("L126").Value = 100000
$ws.Range("N126").Value = -109880
# Row 137
$ws.Range("H137").Value = 13896
$ws.Range("I137").Value = 11349.5
$ws.Range("K137").Value = 34048.5
$ws.Range("M137").Value = -31498.5
# Row 138
$ws.Range("H138").Value = 3193.5366
$ws.Range("I138").Value = 2272.8125
$ws.Range("J138").Value = 3782.8
$ws.Range("K138").Value = 6818.4375
$ws.Range("L138").Value = 11348.4
$ws.Range("M138").Value = -1678.4375
$ws.Range("N138").Value = -21628.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4686.2256
$ws.Range("I61").Value = 4206
$ws.Range("J61").Value = 7183.4
$ws.Range("K61").Value = 4206
$ws.Range("L61").Value = 7183.4
$ws.Range("M61").Value = -3994
$ws.Range("N61").Value = -7607.4
# Row 74
$ws.Range("H74").Value = 201463.48
$ws.Range("I74").Value = 306777.44
$ws.Range("J74").Value = 8387.944
$ws.Range("K74").Value = 306777.44
$ws.Range("L74").Value = 8387.944
$ws.Range("M74").Value = -305903.44
$ws.Range("N74").Value = -10135.944
# Row 77
$ws.Range("H77").Value = 201463.48
$ws.Range("I77").Value = 306777.44
$ws.Range("J77").Value = 8387.944
$ws.Range("K77").Value = 1533887.2
$ws.Range("L77").Value = 41939.72
$ws.Range("M77").Value = -1529519.2
$ws.Range("N77").Value = -50675.72
# Row 92
$ws.Range("H92").Value = 34999.5
$ws.Range("J92").Value = 34999.5
$ws.Range("L92").Value = 34999.5
$ws.Range("N92").Value = -39991.5
# Row 136
$ws.Range("H136").Value = 4686.2256
$ws.Range("I136").Value = 4206
$ws.Range("J136").Value = 7183.4
$ws.Range("K136").Value = 12618
$ws.Range("L136").Value = 21550.2
$ws.Range("M136").Value = -10068
$ws.Range("N136").Value = -26650.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3010.2
$ws.Range("I20").Value = 3429
$ws.Range("K20").Value = 3429
$ws.Range("M20").Value = -3182
# Row 86
$ws.Range("H86").Value = 201488.8
$ws.Range("I86").Value = 1654.2222
$ws.Range("K86").Value = 1654.2222
$ws.Range("M86").Value = -531.2221999999999
# Row 89
$ws.Range("H89").Value = 201488.8
$ws.Range("I89").Value = 1654.2222
$ws.Range("K89").Value = 8271.110999999999
$ws.Range("M89").Value = -2655.110999999999
# Row 102
$ws.Range("H102").Value = 27248.75
$ws.Range("I102").Value = 12998.333
$ws.Range("K102").Value = 12998.333
$ws.Range("M102").Value = -9753.333000000001
# Row 107
$ws.Range("H107").Value = 1346.0741
$ws.Range("I107").Value = 1266.2632
$ws.Range("K107").Value = 1266.2632
$ws.Range("M107").Value = 653.7367999999999
# Row 134
$ws.Range("H134").Value = 4014.6177
$ws.Range("I134").Value = 841
$ws.Range("K134").Value = 2523
$ws.Range("M134").Value = 12

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5383.3335
$ws.Range("I58").Value = 3735.2856
$ws.Range("J58").Value = 6432.091
$ws.Range("K58").Value = 3735.2856
$ws.Range("L58").Value = 6432.091
$ws.Range("M58").Value = -3532.2856
$ws.Range("N58").Value = -6838.091
# Row 125
$ws.Range("H125").Value = 100679.29
$ws.Range("J125").Value = 100679.29
$ws.Range("L125").Value = 100679.29
$ws.Range("N125").Value = -105599.29
# Row 132
$ws.Range("H132").Value = 216117.86
$ws.Range("I132").Value = 150762.5
$ws.Range("J132").Value = 608250
$ws.Range("K132").Value = 452287.5
$ws.Range("L132").Value = 1824750
$ws.Range("M132").Value = -449757.5
$ws.Range("N132").Value = -1829810
# Row 134
$ws.Range("H134").Value = 3681.9412
$ws.Range("I134").Value = 2723.6
$ws.Range("J134").Value = 6344
$ws.Range("K134").Value = 8170.799999999999
$ws.Range("L134").Value = 19032
$ws.Range("M134").Value = -5635.799999999999
$ws.Range("N134").Value = -24102
# Row 136
$ws.Range("H136").Value = 5383.3335
$ws.Range("I136").Value = 3735.2856
$ws.Range("J136").Value = 6432.091
$ws.Range("K136").Value = 11205.8568
$ws.Range("L136").Value = 19296.273
$ws.Range("M136").Value = -8655.856800000001
$ws.Range("N136").Value = -24396.273

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 100
$ws.Range("H100").Value = 499
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 499
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 1497
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -3119
# Row 122
$ws.Range("H122").Value = 17773.75
$ws.Range("I122").Value = 34860
$ws.Range("K122").Value = 313740
$ws.Range("M122").Value = -311290

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 12985.75
$ws.Range("I99").Value = 8425.444
$ws.Range("J99").Value = 26666.666
$ws.Range("K99").Value = 8425.444
$ws.Range("L99").Value = 26666.666
$ws.Range("M99").Value = -6179.444
$ws.Range("N99").Value = -31158.666
# Row 102
$ws.Range("H102").Value = 2392.5557
$ws.Range("J102").Value = 2175.6667
$ws.Range("L102").Value = 2175.6667
$ws.Range("N102").Value = -5419.6667
# Row 122
$ws.Range("H122").Value = 10163.488
$ws.Range("I122").Value = 10662.719
$ws.Range("K122").Value = 31988.157
$ws.Range("M122").Value = -29538.157
# Row 126
$ws.Range("H126").Value = 2949.6
$ws.Range("I126").Value = 2238.3
$ws.Range("J126").Value = 3660.9
$ws.Range("K126").Value = 6714.900000000001
$ws.Range("L126").Value = 10982.7
$ws.Range("M126").Value = -4244.900000000001
$ws.Range("N126").Value = -15922.7
# Row 130
$ws.Range("H130").Value = 81250
$ws.Range("J130").Value = 81250
$ws.Range("L130").Value = 81250
$ws.Range("N130").Value = -91290

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4071.625
$ws.Range("I7").Value = 3824
$ws.Range("K7").Value = 3824
$ws.Range("M7").Value = -3712
# Row 16
$ws.Range("H16").Value = 1618.3125
$ws.Range("I16").Value = 1618.3125
$ws.Range("K16").Value = 1618.3125
$ws.Range("M16").Value = -1448.3125
# Row 55
$ws.Range("H55").Value = 908.3077
$ws.Range("I55").Value = 1027.8
$ws.Range("J55").Value = 510
$ws.Range("K55").Value = 1027.8
$ws.Range("L55").Value = 510
$ws.Range("M55").Value = -854.8
$ws.Range("N55").Value = -856
# Row 68
$ws.Range("H68").Value = 3499.4
$ws.Range("J68").Value = 3499
$ws.Range("L68").Value = 3499
$ws.Range("N68").Value = -4997
# Row 71
$ws.Range("H71").Value = 3499.4
$ws.Range("J71").Value = 3499
$ws.Range("L71").Value = 17495
$ws.Range("N71").Value = -24983
# Row 122
$ws.Range("H122").Value = 8152.1
$ws.Range("I122").Value = 8152.1
$ws.Range("K122").Value = 24456.3
$ws.Range("M122").Value = -22006.3
# Row 126
$ws.Range("H126").Value = 4071.625
$ws.Range("I126").Value = 3824
$ws.Range("K126").Value = 11472
$ws.Range("M126").Value = -9002
# Row 132
$ws.Range("H132").Value = 5147.3447
$ws.Range("I132").Value = 3474.96
$ws.Range("K132").Value = 10424.88
$ws.Range("M132").Value = -7894.880000000001
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 3628.6316
$ws.Range("I126").Value = 3246.5
$ws.Range("J126").Value = 5666.6665
$ws.Range("K126").Value = 9739.5
$ws.Range("L126").Value = 16999.9995
$ws.Range("M126").Value = -7269.5
$ws.Range("N126").Value = -21939.9995
# Row 132
$ws.Range("H132").Value = 3171.8044
$ws.Range("I132").Value = 2000.1111
$ws.Range("J132").Value = 7389.9
$ws.Range("K132").Value = 6000.3333
$ws.Range("L132").Value = 22169.7
$ws.Range("M132").Value = -3470.3333
$ws.Range("N132").Value = -27229.7
# Row 136
$ws.Range("H136").Value = 1980.1945
$ws.Range("I136").Value = 796.1
$ws.Range("J136").Value = 7900.6665
$ws.Range("K136").Value = 2388.3
$ws.Range("L136").Value = 23701.9995
$ws.Range("M136").Value = 161.6999999999998
$ws.Range("N136").Value = -28801.9995
